# Insert a new data row at row 65 (pushes existing rows 65-150 down to 66-151)
# and populate it with a new "Ajo" (garlic) price record for
# "Feria Lagunitas de Puerto Montt", matching the weekly update described
# in the commit message ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(65).Insert()

$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44467
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112003
$ws.Range("G65").Value = "Ajo"
$ws.Range("H65").Value = "Chino"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 250
$ws.Range("K65").Value = 17000
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = 17000
$ws.Range("N65").Value = "$/caja 10 kilos"
$ws.Range("O65").Value = "China"
$ws.Range("P65").Value = 1700
$ws.Range("Q65").Value = 10
$ws.Range("R65").Value = "Hortaliza"
